# Auto-generated Excel COM-interop edit script
# Applies the numeric updates to the leve-crafting profit tables across all 8 sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 3483.3333
$ws.Range("I131").Value = 850
$ws.Range("J131").Value = 4800
$ws.Range("K131").Value = 2550
$ws.Range("L131").Value = 14400
$ws.Range("M131").Value = 2490
$ws.Range("N131").Value = -24480

$ws.Range("H138").Value = 2277.7122
$ws.Range("I138").Value = 1364.3667
$ws.Range("J138").Value = 3038.8333
$ws.Range("K138").Value = 4093.1001
$ws.Range("L138").Value = 9116.499899999999
$ws.Range("M138").Value = 1046.8999
$ws.Range("N138").Value = -19396.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20193.64
$ws.Range("I32").Value = 3492.836
$ws.Range("J32").Value = 121916.73
$ws.Range("K32").Value = 3492.836
$ws.Range("L32").Value = 121916.73
$ws.Range("M32").Value = -3205.836
$ws.Range("N32").Value = -122490.73

$ws.Range("H110").Value = 2244
$ws.Range("I110").Value = 1348.8823
$ws.Range("J110").Value = 4417.857
$ws.Range("K110").Value = 1348.8823
$ws.Range("L110").Value = 4417.857
$ws.Range("M110").Value = 696.1177
$ws.Range("N110").Value = -8507.857

$ws.Range("H122").Value = 1770.5
$ws.Range("I122").Value = 1713.125
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5139.375
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2689.375
$ws.Range("N122").Value = -10900

$ws.Range("H132").Value = 1444.5209
$ws.Range("I132").Value = 1004.9286
$ws.Range("J132").Value = 2059.95
$ws.Range("K132").Value = 3014.7858
$ws.Range("L132").Value = 6179.849999999999
$ws.Range("M132").Value = -484.7857999999997
$ws.Range("N132").Value = -11239.85

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1822.3793
$ws.Range("I99").Value = 1175.7368
$ws.Range("K99").Value = 1175.7368
$ws.Range("M99").Value = 322.2632000000001

$ws.Range("H107").Value = 1536.5151
$ws.Range("I107").Value = 1516.2084
$ws.Range("J107").Value = 1590.6666
$ws.Range("K107").Value = 1516.2084
$ws.Range("L107").Value = 1590.6666
$ws.Range("M107").Value = 403.7916
$ws.Range("N107").Value = -5430.6666

$ws.Range("H134").Value = 1533.3226
$ws.Range("I134").Value = 1300.7916
$ws.Range("J134").Value = 2330.5715
$ws.Range("K134").Value = 3902.3748
$ws.Range("L134").Value = 6991.7145
$ws.Range("M134").Value = -1367.3748
$ws.Range("N134").Value = -12061.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 16046.154
$ws.Range("I17").Value = 455.55554
$ws.Range("J17").Value = 51125
$ws.Range("K17").Value = 1366.66662
$ws.Range("L17").Value = 153375
$ws.Range("M17").Value = -1197.66662
$ws.Range("N17").Value = -153713

$ws.Range("H21").Value = 2750
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 2750
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 8250
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -8596

$ws.Range("H69").Value = 4637
$ws.Range("J69").Value = 5228
$ws.Range("L69").Value = 15684
$ws.Range("N69").Value = -17306

$ws.Range("H72").Value = 4637
$ws.Range("J72").Value = 5228
$ws.Range("L72").Value = 47052
$ws.Range("N72").Value = -55164

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 1000
$ws.Range("I21").Value = 1000
$ws.Range("K21").Value = 1000
$ws.Range("M21").Value = -827

$ws.Range("H30").Value = 1000
$ws.Range("I30").Value = 1000
$ws.Range("K30").Value = 1000
$ws.Range("M30").Value = -895

$ws.Range("H70").Value = 27682.564
$ws.Range("I70").Value = 40862
$ws.Range("J70").Value = 5200
$ws.Range("K70").Value = 40862
$ws.Range("L70").Value = 5200
$ws.Range("M70").Value = -40592
$ws.Range("N70").Value = -5740

$ws.Range("H73").Value = 27682.564
$ws.Range("I73").Value = 40862
$ws.Range("J73").Value = 5200
$ws.Range("K73").Value = 40862
$ws.Range("L73").Value = 5200
$ws.Range("M73").Value = -39926
$ws.Range("N73").Value = -7072

$ws.Range("H80").Value = 2922.2856
$ws.Range("I80").Value = 2655.9
$ws.Range("J80").Value = 3588.25
$ws.Range("K80").Value = 2655.9
$ws.Range("L80").Value = 3588.25
$ws.Range("M80").Value = -1657.9
$ws.Range("N80").Value = -5584.25

$ws.Range("H83").Value = 2922.2856
$ws.Range("I83").Value = 2655.9
$ws.Range("J83").Value = 3588.25
$ws.Range("K83").Value = 13279.5
$ws.Range("L83").Value = 17941.25
$ws.Range("M83").Value = -8287.5
$ws.Range("N83").Value = -27925.25

$ws.Range("H107").Value = 545.8214
$ws.Range("I107").Value = 445.53333
$ws.Range("J107").Value = 661.53845
$ws.Range("K107").Value = 445.53333
$ws.Range("L107").Value = 661.53845
$ws.Range("M107").Value = 1474.46667
$ws.Range("N107").Value = -4501.53845

$ws.Range("H132").Value = 2816.4614
$ws.Range("I132").Value = 2263.647
$ws.Range("J132").Value = 3860.6667
$ws.Range("K132").Value = 6790.941
$ws.Range("L132").Value = 11582.0001
$ws.Range("M132").Value = -4260.941
$ws.Range("N132").Value = -16642.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 20887
$ws.Range("I30").Value = 20887
$ws.Range("K30").Value = 20887
$ws.Range("M30").Value = -20779

$ws.Range("H55").Value = 748.4
$ws.Range("I55").Value = 298.4
$ws.Range("J55").Value = 973.4
$ws.Range("K55").Value = 298.4
$ws.Range("L55").Value = 973.4
$ws.Range("M55").Value = -125.4
$ws.Range("N55").Value = -1319.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

$ws.Range("H80").Value = 26725
$ws.Range("J80").Value = 26725
$ws.Range("L80").Value = 26725
$ws.Range("N80").Value = -28721

$ws.Range("H83").Value = 26725
$ws.Range("J83").Value = 26725
$ws.Range("L83").Value = 80175
$ws.Range("N83").Value = -90159
